$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1755458.1
$ws.Range("I19").Value = 2924840.8
$ws.Range("J19").Value = 1384
$ws.Range("K19").Value = 2924840.8
$ws.Range("L19").Value = 1384
$ws.Range("M19").Value = -2924665.8
$ws.Range("N19").Value = -1734
$ws.Range("H74").Value = 4766348
$ws.Range("I74").Value = 12503538
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 12503538
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -12502602
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4766348
$ws.Range("I77").Value = 12503538
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 62517690
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -62513010
$ws.Range("N77").Value = -34360
$ws.Range("H116").Value = 599094.75
$ws.Range("I116").Value = 3337601.8
$ws.Range("J116").Value = 12271.857
$ws.Range("K116").Value = 3337601.8
$ws.Range("L116").Value = 12271.857
$ws.Range("M116").Value = -3334159.8
$ws.Range("N116").Value = -19155.857
$ws.Range("H138").Value = 2292.1702
$ws.Range("I138").Value = 771.875
$ws.Range("J138").Value = 3076.8386
$ws.Range("K138").Value = 2315.625
$ws.Range("L138").Value = 9230.515800000001
$ws.Range("M138").Value = 2824.375
$ws.Range("N138").Value = -19510.5158

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10336.667
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 10336.667
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = 10336.667
$ws.Range("N3").Value = -10566.667
$ws.Range("H61").Value = 1140.48
$ws.Range("I61").Value = 853.05554
$ws.Range("K61").Value = 853.05554
$ws.Range("M61").Value = -641.05554
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H132").Value = 3820.087
$ws.Range("I132").Value = 2216.5
$ws.Range("J132").Value = 5053.615
$ws.Range("K132").Value = 6649.5
$ws.Range("L132").Value = 15160.845
$ws.Range("M132").Value = -4119.5
$ws.Range("N132").Value = -20220.845
$ws.Range("H136").Value = 1140.48
$ws.Range("I136").Value = 853.05554
$ws.Range("K136").Value = 2559.16662
$ws.Range("M136").Value = -9.166619999999966

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25364.867
$ws.Range("I82").Value = 4571
$ws.Range("J82").Value = 30563.334
$ws.Range("K82").Value = 4571
$ws.Range("L82").Value = 30563.334
$ws.Range("M82").Value = -4188
$ws.Range("N82").Value = -31329.334
$ws.Range("H85").Value = 25364.867
$ws.Range("I85").Value = 4571
$ws.Range("J85").Value = 30563.334
$ws.Range("K85").Value = 4571
$ws.Range("L85").Value = 30563.334
$ws.Range("M85").Value = -3245
$ws.Range("N85").Value = -33215.334
$ws.Range("H86").Value = 1650
$ws.Range("I86").Value = 1580
$ws.Range("K86").Value = 1580
$ws.Range("M86").Value = -457
$ws.Range("H89").Value = 1650
$ws.Range("I89").Value = 1580
$ws.Range("K89").Value = 7900
$ws.Range("M89").Value = -2284
$ws.Range("H134").Value = 1713.9342
$ws.Range("I134").Value = 1049.25
$ws.Range("J134").Value = 3154.0833
$ws.Range("K134").Value = 3147.75
$ws.Range("L134").Value = 9462.249899999999
$ws.Range("M134").Value = -612.75
$ws.Range("N134").Value = -14532.2499

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 537.19354
$ws.Range("I22").Value = 340.18182
$ws.Range("K22").Value = 340.18182
$ws.Range("M22").Value = 9.818179999999984
$ws.Range("H31").Value = 7145276.5
$ws.Range("I31").Value = 1321.341
$ws.Range("K31").Value = 1321.341
$ws.Range("M31").Value = -1026.341
$ws.Range("H34").Value = 7145276.5
$ws.Range("I34").Value = 1321.341
$ws.Range("K34").Value = 1321.341
$ws.Range("M34").Value = -1119.341
$ws.Range("H134").Value = 3180.5344
$ws.Range("I134").Value = 4230.8
$ws.Range("J134").Value = 2055.25
$ws.Range("K134").Value = 12692.4
$ws.Range("L134").Value = 6165.75
$ws.Range("M134").Value = -10157.4
$ws.Range("N134").Value = -11235.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1363.8422
$ws.Range("I5").Value = 327.3158
$ws.Range("J5").Value = 2400.3684
$ws.Range("K5").Value = 981.9474
$ws.Range("L5").Value = 7201.1052
$ws.Range("M5").Value = -869.9474
$ws.Range("N5").Value = -7425.1052
$ws.Range("H131").Value = 858.52
$ws.Range("I131").Value = 560.36365
$ws.Range("J131").Value = 895.3708
$ws.Range("K131").Value = 1681.09095
$ws.Range("L131").Value = 2686.1124
$ws.Range("M131").Value = 3358.90905
$ws.Range("N131").Value = -12766.1124
$ws.Range("H132").Value = 1656.9474
$ws.Range("I132").Value = 708.3
$ws.Range("J132").Value = 2711
$ws.Range("K132").Value = 6374.7
$ws.Range("L132").Value = 24399
$ws.Range("M132").Value = -3844.7
$ws.Range("N132").Value = -29459
$ws.Range("H135").Value = 1363.8422
$ws.Range("I135").Value = 327.3158
$ws.Range("J135").Value = 2400.3684
$ws.Range("K135").Value = 2945.8422
$ws.Range("L135").Value = 21603.3156
$ws.Range("M135").Value = -410.8422
$ws.Range("N135").Value = -26673.3156

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 29920.625
$ws.Range("J46").Value = 29920.625
$ws.Range("L46").Value = 29920.625
$ws.Range("N46").Value = -30232.625
$ws.Range("H70").Value = 6393.9487
$ws.Range("I70").Value = 5811.4443
$ws.Range("K70").Value = 5811.4443
$ws.Range("M70").Value = -5541.4443
$ws.Range("H73").Value = 6393.9487
$ws.Range("I73").Value = 5811.4443
$ws.Range("K73").Value = 5811.4443
$ws.Range("M73").Value = -4875.4443
$ws.Range("H126").Value = 3253.4343
$ws.Range("I126").Value = 3000.15
$ws.Range("J126").Value = 4319.8945
$ws.Range("K126").Value = 9000.450000000001
$ws.Range("L126").Value = 12959.6835
$ws.Range("M126").Value = -6530.450000000001
$ws.Range("N126").Value = -17899.6835
$ws.Range("H132").Value = 6974
$ws.Range("I132").Value = 5666
$ws.Range("J132").Value = 7235.6
$ws.Range("K132").Value = 16998
$ws.Range("L132").Value = 21706.8
$ws.Range("M132").Value = -14468
$ws.Range("N132").Value = -26766.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9997.159
$ws.Range("I132").Value = 10824.885
$ws.Range("J132").Value = 8801.556
$ws.Range("K132").Value = 32474.655
$ws.Range("L132").Value = 26404.668
$ws.Range("M132").Value = -29944.655
$ws.Range("N132").Value = -31464.668
$ws.Range("H136").Value = 2191.132
$ws.Range("I136").Value = 1268.1395
$ws.Range("K136").Value = 3804.4185
$ws.Range("M136").Value = -1254.4185

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1867.194
$ws.Range("I136").Value = 675.3333
$ws.Range("K136").Value = 2025.9999
$ws.Range("M136").Value = 524.0001

Write-Host "Applied changes: set=175 delete=1 add=2"